$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 70
$ws.Cells.Item(3, 6).Value = 229
$ws.Cells.Item(4, 6).Value = 68
$ws.Cells.Item(5, 6).Value = 9107
$ws.Cells.Item(6, 6).Value = 534
$ws.Cells.Item(7, 6).Value = 91
$ws.Cells.Item(9, 6).Value = 212
$ws.Cells.Item(10, 6).Value = 328
$ws.Cells.Item(11, 6).Value = 378
$ws.Cells.Item(14, 6).Value = 411
$ws.Cells.Item(15, 6).Value = 11651
$ws.Cells.Item(16, 6).Value = 11651
$ws.Cells.Item(20, 6).Value = 82
$ws.Cells.Item(21, 6).Value = 30
$ws.Cells.Item(24, 6).Value = 150
$ws.Cells.Item(26, 6).Value = 215
$ws.Cells.Item(30, 6).Value = 2704
$ws.Cells.Item(33, 6).Value = 2090
$ws.Cells.Item(34, 6).Value = 56
$ws.Cells.Item(36, 6).Value = 2129
$ws.Cells.Item(37, 6).Value = 956
$ws.Cells.Item(38, 6).Value = 4161
$ws.Cells.Item(39, 6).Value = 317
$ws.Cells.Item(40, 6).Value = 3046
$ws.Cells.Item(41, 6).Value = 1285
$ws.Cells.Item(43, 6).Value = 86
$ws.Cells.Item(44, 6).Value = 383
$ws.Cells.Item(45, 6).Value = 438
$ws.Cells.Item(48, 6).Value = 175
$ws.Cells.Item(50, 6).Value = 113

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 19
$ws.Cells.Item(6, 6).Value = 12
$ws.Cells.Item(18, 6).Value = 4
$ws.Cells.Item(20, 6).Value = 70
$ws.Cells.Item(22, 6).Value = 29

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 42

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 70
$ws.Cells.Item(3, 6).Value = 19
$ws.Cells.Item(5, 6).Value = 12
$ws.Cells.Item(7, 6).Value = 229
$ws.Cells.Item(8, 6).Value = 68
$ws.Cells.Item(9, 6).Value = 9107
$ws.Cells.Item(10, 6).Value = 534
$ws.Cells.Item(12, 6).Value = 91
$ws.Cells.Item(13, 6).Value = 212
$ws.Cells.Item(14, 6).Value = 328
$ws.Cells.Item(15, 6).Value = 378
$ws.Cells.Item(18, 6).Value = 11651
$ws.Cells.Item(22, 6).Value = 82
$ws.Cells.Item(23, 6).Value = 30
$ws.Cells.Item(26, 6).Value = 150
$ws.Cells.Item(28, 6).Value = 215
$ws.Cells.Item(34, 6).Value = 2090
$ws.Cells.Item(35, 6).Value = 56
$ws.Cells.Item(37, 6).Value = 2129
$ws.Cells.Item(38, 6).Value = 956
$ws.Cells.Item(39, 6).Value = 4
$ws.Cells.Item(40, 6).Value = 4161
$ws.Cells.Item(41, 6).Value = 317
$ws.Cells.Item(42, 6).Value = 3046
$ws.Cells.Item(43, 6).Value = 1285
$ws.Cells.Item(44, 6).Value = 86
$ws.Cells.Item(45, 6).Value = 383
$ws.Cells.Item(48, 6).Value = 175
$ws.Cells.Item(50, 6).Value = 113
